# "updated datasets for rev 2"
# The experiment previously labelled "dur1" (rows 2-15, column C / "experiment")
# is renamed to "pd" on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2:C15").Value = "pd"
